$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 23.86000000000029
$ws.Range("H2").Value = 0.001027943970288381
$ws.Range("I2").Value = 0.001027943970288381
$ws.Range("L2").Value = 36.28051904447643
$ws.Range("M2").Value = "[15.302038547125854, 57.258999541827]"
$ws.Range("N2").Value = 0.001115128605711613
$ws.Range("O2").Value = 0.001115128605711613
$ws.Range("P2").Value = 1.515763422452733
$ws.Range("Q2").Value = "[0.748447499053424, 2.283079345852042]"
$ws.Range("R2").Value = 0.000249252601888017
$ws.Range("S2").Value = 0.000249252601888017
$ws.Range("T2").Value = 49.4355076277879
$ws.Range("U2").Value = "[36.601207986566244, 62.269807269009554]"
$ws.Range("V2").Value = [double]"7.774882959665774e-10"
$ws.Range("W2").Value = [double]"7.774882959665774e-10"
$ws.Range("X2").Value = 18.1039839839842
$ws.Range("Y2").Value = 15.19015015015033
$ws.Range("Z2").Value = 21.01781781781807

# Row 3
$ws.Range("F3").Value = 23.86000000000029
$ws.Range("H3").Value = 0.02813418718459448
$ws.Range("I3").Value = 0.02813418718459448
$ws.Range("L3").Value = 27.16720309037451
$ws.Range("M3").Value = "[3.319335637886411, 51.0150705428626]"
$ws.Range("N3").Value = 0.02648432784939914
$ws.Range("O3").Value = 0.02648432784939914
$ws.Range("P3").Value = 1.440289725069194
$ws.Range("Q3").Value = "[0.19497371824080645, 2.6856057318975814]"
$ws.Range("R3").Value = 0.02438267407129469
$ws.Range("S3").Value = 0.02438267407129469
$ws.Range("T3").Value = 48.9097788443317
$ws.Range("U3").Value = "[35.04148838012806, 62.77806930853535]"
$ws.Range("V3").Value = [double]"7.165113169449455e-09"
$ws.Range("W3").Value = [double]"7.165113169449455e-09"
$ws.Range("X3").Value = 18.39059059059081
$ws.Range("Y3").Value = 13.66158158158174
$ws.Range("Z3").Value = 23.11959959959989

# Row 4
$ws.Range("F4").Value = 23.86000000000029
$ws.Range("H4").Value = 0.0001275412456203107
$ws.Range("I4").Value = 0.0001275412456203107
$ws.Range("L4").Value = 41.82691843955279
$ws.Range("M4").Value = "[19.673872142762356, 63.97996473634323]"
$ws.Range("N4").Value = 0.0004282886995232715
$ws.Range("O4").Value = 0.0004282886995232715
$ws.Range("P4").Value = 1.742184514603348
$ws.Range("Q4").Value = "[1.1258159859711165, 2.35855304323558]"
$ws.Range("R4").Value = [double]"8.899715004950792e-07"
$ws.Range("S4").Value = [double]"8.899715004950792e-07"
$ws.Range("T4").Value = 56.15631341494793
$ws.Range("U4").Value = "[43.51497784889994, 68.79764898099592]"
$ws.Range("V4").Value = [double]"1.50464085635349e-11"
$ws.Range("W4").Value = [double]"1.50464085635349e-11"
$ws.Range("X4").Value = 17.24416416416437
$ws.Range("Y4").Value = 14.90354354354372
$ws.Range("Z4").Value = 19.58478478478502

# Row 5
$ws.Range("F5").Value = 23.86000000000029
$ws.Range("H5").Value = 0.01168834708232536
$ws.Range("I5").Value = 0.01168834708232536
$ws.Range("L5").Value = 29.85635109074131
$ws.Range("M5").Value = "[5.721811126121381, 53.990891055361246]"
$ws.Range("N5").Value = 0.01646962308390165
$ws.Range("O5").Value = 0.01646962308390165
$ws.Range("P5").Value = 1.855395060678656
$ws.Range("Q5").Value = "[0.8365001460008861, 2.874289975356425]"
$ws.Range("R5").Value = 0.0006447774461866285
$ws.Range("S5").Value = 0.0006447774461866285
$ws.Range("T5").Value = 57.56097141035458
$ws.Range("U5").Value = "[44.048489066390346, 71.07345375431882]"
$ws.Range("V5").Value = [double]"5.016520532308277e-11"
$ws.Range("W5").Value = [double]"5.016520532308277e-11"
$ws.Range("X5").Value = 16.81425425425446
$ws.Range("Y5").Value = 12.94506506506523
$ws.Range("Z5").Value = 20.6834434434437

# Row 6
$ws.Range("F6").Value = 23.86000000000029
$ws.Range("H6").Value = [double]"3.278643245119284e-06"
$ws.Range("I6").Value = [double]"3.278643245119284e-06"
$ws.Range("L6").Value = 48.7677008320089
$ws.Range("M6").Value = "[29.605350705420904, 67.9300509585969]"
$ws.Range("N6").Value = [double]"6.043353582896316e-06"
$ws.Range("O6").Value = [double]"6.043353582896316e-06"
$ws.Range("P6").Value = 1.66671081721981
$ws.Range("Q6").Value = "[1.2012896833546556, 2.132131951084965]"
$ws.Range("R6").Value = [double]"4.934422426217111e-09"
$ws.Range("S6").Value = [double]"4.934422426217111e-09"
$ws.Range("T6").Value = 59.59679294980808
$ws.Range("U6").Value = "[47.70267917161412, 71.49090672800203]"
$ws.Range("V6").Value = [double]"3.9013237085328e-13"
$ws.Range("W6").Value = [double]"3.9013237085328e-13"
$ws.Range("X6").Value = 17.53077077077099
$ws.Range("Y6").Value = 15.76336336336356
$ws.Range("Z6").Value = 19.29817817817841

# Row 7
$ws.Range("F7").Value = 23.86000000000029
$ws.Range("H7").Value = 0.03099943804893068
$ws.Range("I7").Value = 0.03099943804893068
$ws.Range("L7").Value = 25.30270668385357
$ws.Range("M7").Value = "[2.2034458488139634, 48.40196751889317]"
$ws.Range("N7").Value = 0.03251486979480278
$ws.Range("O7").Value = 0.03251486979480278
$ws.Range("P7").Value = 1.767342413731195
$ws.Range("Q7").Value = "[0.5471843060306547, 2.987500521431736]"
$ws.Range("R7").Value = 0.005491244830893249
$ws.Range("S7").Value = 0.005491244830893249
$ws.Range("T7").Value = 53.22967857153009
$ws.Range("U7").Value = "[40.12954357039311, 66.32981357266706]"
$ws.Range("V7").Value = [double]"1.864139953511312e-10"
$ws.Range("W7").Value = [double]"1.864139953511312e-10"
$ws.Range("X7").Value = 17.14862862862883
$ws.Range("Y7").Value = 12.5151551551553
$ws.Range("Z7").Value = 21.78210210210236

# Row 8
$ws.Range("F8").Value = 23.86000000000029
$ws.Range("H8").Value = 0.001138975248351271
$ws.Range("I8").Value = 0.001138975248351271
$ws.Range("L8").Value = 46.80024820076011
$ws.Range("M8").Value = "[15.76762588172474, 77.83287051979548]"
$ws.Range("N8").Value = 0.003961676057310903
$ws.Range("O8").Value = 0.003961676057310903
$ws.Range("P8").Value = 1.905710858934349
$ws.Range("Q8").Value = "[1.1887107337907334, 2.6227109840779654]"
$ws.Range("R8").Value = [double]"2.813854661276594e-06"
$ws.Range("S8").Value = [double]"2.813854661276594e-06"
$ws.Range("T8").Value = 63.0153576098115
$ws.Range("U8").Value = "[46.303269721149135, 79.72744549847386]"
$ws.Range("V8").Value = [double]"1.350352718532122e-09"
$ws.Range("W8").Value = [double]"1.350352718532122e-09"
$ws.Range("X8").Value = 16.62318318318339
$ws.Range("Y8").Value = 13.90042042042059
$ws.Range("Z8").Value = 19.34594594594618

# Row 9
$ws.Range("F9").Value = 23.86000000000029
$ws.Range("H9").Value = 0.008960261999606978
$ws.Range("I9").Value = 0.008960261999606978
$ws.Range("L9").Value = 33.35171066380452
$ws.Range("M9").Value = "[7.780198870568384, 58.923222457040666]"
$ws.Range("N9").Value = 0.01173833349337294
$ws.Range("O9").Value = 0.01173833349337294
$ws.Range("P9").Value = 1.415131825941349
$ws.Range("Q9").Value = "[0.42139481039142535, 2.4088688414912722]"
$ws.Range("R9").Value = 0.006263414616181961
$ws.Range("S9").Value = 0.006263414616181961
$ws.Range("T9").Value = 54.85779905250182
$ws.Range("U9").Value = "[40.24542084069817, 69.47017726430546]"
$ws.Range("V9").Value = [double]"1.510473746080265e-09"
$ws.Range("W9").Value = [double]"1.510473746080265e-09"
$ws.Range("X9").Value = 18.48612612612635
$ws.Range("Y9").Value = 14.71247247247265
$ws.Range("Z9").Value = 22.25977977978005

# Row 10
$ws.Range("F10").Value = 23.86000000000029
$ws.Range("H10").Value = 0.007761162267299615
$ws.Range("I10").Value = 0.007761162267299615
$ws.Range("L10").Value = 31.13836939445484
$ws.Range("M10").Value = "[5.197068754063459, 57.079670034846224]"
$ws.Range("N10").Value = 0.01973570625809096
$ws.Range("O10").Value = 0.01973570625809096
$ws.Range("P10").Value = 2.018921405009656
$ws.Range("Q10").Value = "[1.1509738850989635, 2.886868924920348]"
$ws.Range("R10").Value = [double]"2.605158405022401e-05"
$ws.Range("S10").Value = [double]"2.605158405022401e-05"
$ws.Range("T10").Value = 58.77863108267499
$ws.Range("U10").Value = "[45.349584527740745, 72.20767763760924]"
$ws.Range("V10").Value = [double]"2.311506541730068e-11"
$ws.Range("W10").Value = [double]"2.311506541730068e-11"
$ws.Range("X10").Value = 16.19327327327348
$ws.Range("Y10").Value = 12.89729729729746
$ws.Range("Z10").Value = 19.48924924924949

# Row 11
$ws.Range("F11").Value = 24.05000000000032
$ws.Range("H11").Value = 0.004743233935778135
$ws.Range("I11").Value = 0.004743233935778135
$ws.Range("L11").Value = 31.44761830058168
$ws.Range("M11").Value = "[8.089150588709032, 54.80608601245432]"
$ws.Range("N11").Value = 0.009448071551219961
$ws.Range("O11").Value = 0.009448071551219961
$ws.Range("P11").Value = 2.031500354573581
$ws.Range("Q11").Value = "[1.2138686329185795, 2.8491320762285817]"
$ws.Range("R11").Value = [double]"9.069152162810568e-06"
$ws.Range("S11").Value = [double]"9.069152162810568e-06"
$ws.Range("T11").Value = 51.85024147986725
$ws.Range("U11").Value = "[38.957435444296074, 64.74304751543843]"
$ws.Range("V11").Value = [double]"2.466653548083286e-10"
$ws.Range("W11").Value = [double]"2.466653548083286e-10"
$ws.Range("X11").Value = 16.27407407407429
$ws.Range("Y11").Value = 13.14444444444462
$ws.Range("Z11").Value = 19.40370370370396

# Row 12
$ws.Range("F12").Value = 24.05000000000032
$ws.Range("H12").Value = 0.0005162325536504087
$ws.Range("I12").Value = 0.0005162325536504087
$ws.Range("L12").Value = 40.36051305256965
$ws.Range("M12").Value = "[15.518073436139574, 65.20295266899973]"
$ws.Range("N12").Value = 0.002053454580201342
$ws.Range("O12").Value = 0.002053454580201342
$ws.Range("P12").Value = 2.169868799776734
$ws.Range("Q12").Value = "[1.5535002711445012, 2.7862373284089665]"
$ws.Range("R12").Value = [double]"7.48333395250711e-09"
$ws.Range("S12").Value = [double]"7.48333395250711e-09"
$ws.Range("T12").Value = 50.93694511056864
$ws.Range("U12").Value = "[37.382745652169206, 64.49114456896808]"
$ws.Range("V12").Value = [double]"1.471665234120678e-09"
$ws.Range("W12").Value = [double]"1.471665234120678e-09"
$ws.Range("X12").Value = 15.74444444444465
$ws.Range("Y12").Value = 13.38518518518536
$ws.Range("Z12").Value = 18.10370370370395

# Row 13
$ws.Range("F13").Value = 24.05000000000032
$ws.Range("H13").Value = 0.003136220578632165
$ws.Range("I13").Value = 0.003136220578632165
$ws.Range("L13").Value = 37.26814482250946
$ws.Range("M13").Value = "[10.920551565143633, 63.61573807987528]"
$ws.Range("N13").Value = 0.006592946230671082
$ws.Range("O13").Value = 0.006592946230671082
$ws.Range("P13").Value = 1.817658211986887
$ws.Range("Q13").Value = "[0.9497106920761951, 2.6856057318975797]"
$ws.Range("R13").Value = 0.0001176044800996756
$ws.Range("S13").Value = 0.0001176044800996756
$ws.Range("T13").Value = 74.2152817906464
$ws.Range("U13").Value = "[59.57069291119576, 88.85987067009704]"
$ws.Range("V13").Value = [double]"2.726707748479384e-13"
$ws.Range("W13").Value = [double]"2.726707748479384e-13"
$ws.Range("X13").Value = 17.09259259259282
$ws.Range("Y13").Value = 13.77037037037056
$ws.Range("Z13").Value = 20.41481481481508

# Row 14
$ws.Range("F14").Value = 24.05000000000032
$ws.Range("H14").Value = [double]"1.093281123509815e-05"
$ws.Range("I14").Value = [double]"1.093281123509815e-05"
$ws.Range("L14").Value = 51.49451822603526
$ws.Range("M14").Value = "[25.837791355396448, 77.15124509667407]"
$ws.Range("N14").Value = 0.0002043715423702963
$ws.Range("O14").Value = 0.0002043715423702963
$ws.Range("P14").Value = 2.106974051957119
$ws.Range("Q14").Value = "[1.591237119836272, 2.6227109840779654]"
$ws.Range("R14").Value = [double]"1.607427524419336e-10"
$ws.Range("S14").Value = [double]"1.607427524419336e-10"
$ws.Range("T14").Value = 62.1181509287626
$ws.Range("U14").Value = "[48.67565683812134, 75.56064501940386]"
$ws.Range("V14").Value = [double]"4.691358412856061e-12"
$ws.Range("W14").Value = [double]"4.691358412856061e-12"
$ws.Range("X14").Value = 15.9851851851854
$ws.Range("Y14").Value = 14.0111111111113
$ws.Range("Z14").Value = 17.9592592592595
